$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (120) as a formatting template for the two
# new rows, so the new rows pick up the exact same cell styles (e.g. the
# date style on column A) without registering any new style entries.
$template = $ws.Range("A120:H120")

# ---- Row 121 ----
$template.Copy()
$ws.Range("A121:H121").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(121, 1).Value = 45457.2916666667
$ws.Cells.Item(121, 2).Value = 0
$ws.Cells.Item(121, 3).Value = 3.1800000667572
$ws.Cells.Item(121, 4).Value = 3.1800000667572
$ws.Cells.Item(121, 5).Value = 3.1800000667572
$ws.Cells.Item(121, 6).Value = 3.1800000667572

# Column G holds a numeric-looking string ("3.1800000667572"); assigning it
# directly via .Value would auto-coerce it to a number. Route it through a
# text formula first, then collapse the formula to its cached literal value
# via a self copy/paste-special(values) so it lands as plain shared-string
# text (matching the source file) without touching the style table.
$g121 = $ws.Cells.Item(121, 7)
$g121.Formula = '="3.1800000667572"'
$g121.Copy()
$g121.PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item(121, 8).Value = "AGAIN.MI"

# ---- Row 122 ----
$template.Copy()
$ws.Range("A122:H122").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(122, 1).Value = 45460.6269212963
$ws.Cells.Item(122, 2).Value = 2000
$ws.Cells.Item(122, 3).Value = 3.14000010490417
$ws.Cells.Item(122, 4).Value = 3.11999988555908
$ws.Cells.Item(122, 5).Value = 3.14000010490417
$ws.Cells.Item(122, 6).Value = 3.11999988555908

$g122 = $ws.Cells.Item(122, 7)
$g122.Formula = '="3.11999988555908"'
$g122.Copy()
$g122.PasteSpecial(-4163)  # xlPasteValues

$ws.Cells.Item(122, 8).Value = "AGAIN.MI"
